# Auto-generated edit script applying Hades_Profits.xlsx diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 586.4
$ws.Range("I4").Value = 210
$ws.Range("J4").Value = 1151
$ws.Range("K4").Value = 210
$ws.Range("L4").Value = 1151
$ws.Range("M4").Value = -96
$ws.Range("N4").Value = -1379
$ws.Range("H15").Value = 3762.8867
$ws.Range("I15").Value = 3762.8867
$ws.Range("K15").Value = 11288.6601
$ws.Range("M15").Value = -11119.6601
$ws.Range("H19").Value = 499.5
$ws.Range("I19").Value = 193.25
$ws.Range("J19").Value = 652.625
$ws.Range("K19").Value = 193.25
$ws.Range("L19").Value = 652.625
$ws.Range("M19").Value = -18.25
$ws.Range("N19").Value = -1002.625
$ws.Range("H40").Value = 1526.3572
$ws.Range("I40").Value = 1912.6666
$ws.Range("J40").Value = 1480
$ws.Range("K40").Value = 1912.6666
$ws.Range("L40").Value = 1480
$ws.Range("M40").Value = -1737.6666
$ws.Range("N40").Value = -1830
$ws.Range("H132").Value = 2334649.2
$ws.Range("I132").Value = 1336.3889
$ws.Range("J132").Value = 16334527
$ws.Range("K132").Value = 4009.1667
$ws.Range("L132").Value = 49003581
$ws.Range("M132").Value = -1479.1667
$ws.Range("N132").Value = -49008641
$ws.Range("H133").Value = 30000
$ws.Range("J133").Value = 30000
$ws.Range("L133").Value = 30000
$ws.Range("N133").Value = -40120

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H32").Value = 21066.54
$ws.Range("I32").Value = 20553.299
$ws.Range("J32").Value = 24723.375
$ws.Range("K32").Value = 20553.299
$ws.Range("L32").Value = 24723.375
$ws.Range("M32").Value = -20266.299
$ws.Range("N32").Value = -25297.375
$ws.Range("H61").Value = 24440250
$ws.Range("I61").Value = 31282330
$ws.Range("K61").Value = 31282330
$ws.Range("M61").Value = -31282118
$ws.Range("H74").Value = 4943040.5
$ws.Range("I74").Value = 7383509
$ws.Range("J74").Value = 62104.41
$ws.Range("K74").Value = 7383509
$ws.Range("L74").Value = 62104.41
$ws.Range("M74").Value = -7382635
$ws.Range("N74").Value = -63852.41
$ws.Range("H77").Value = 4943040.5
$ws.Range("I77").Value = 7383509
$ws.Range("J77").Value = 62104.41
$ws.Range("K77").Value = 36917545
$ws.Range("L77").Value = 310522.05
$ws.Range("M77").Value = -36913177
$ws.Range("N77").Value = -319258.05
$ws.Range("H102").Value = 35716030
$ws.Range("I102").Value = 47620544
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 47620544
$ws.Range("L102").Value = 2500
$ws.Range("M102").Value = -47618922
$ws.Range("N102").Value = -5744
$ws.Range("H114").Value = 39699
$ws.Range("J114").Value = 39699
$ws.Range("L114").Value = 39699
$ws.Range("N114").Value = -48377
$ws.Range("H136").Value = 24440250
$ws.Range("I136").Value = 31282330
$ws.Range("K136").Value = 93846990
$ws.Range("M136").Value = -93844440

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H80").Value = 387.5
$ws.Range("I80").Value = 101.333336
$ws.Range("J80").Value = 510.14285
$ws.Range("K80").Value = 101.333336
$ws.Range("L80").Value = 510.14285
$ws.Range("M80").Value = 896.666664
$ws.Range("N80").Value = -2506.14285
$ws.Range("H83").Value = 387.5
$ws.Range("I83").Value = 101.333336
$ws.Range("J83").Value = 510.14285
$ws.Range("K83").Value = 506.66668
$ws.Range("L83").Value = 2550.71425
$ws.Range("M83").Value = 4485.33332
$ws.Range("N83").Value = -12534.71425

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 30135.5
$ws.Range("J64").Value = 30135.5
$ws.Range("L64").Value = 30135.5
$ws.Range("N64").Value = -30631.5
$ws.Range("H67").Value = 30135.5
$ws.Range("J67").Value = 30135.5
$ws.Range("L67").Value = 30135.5
$ws.Range("N67").Value = -31851.5
$ws.Range("H69").Value = 17196
$ws.Range("H72").Value = 17196
$ws.Range("H105").Value = 1975.8
$ws.Range("I105").Value = 1975.8
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1975.8
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -228.8
$ws.Range("N105").ClearContents()
$ws.Range("H132").Value = 22785.188
$ws.Range("I132").Value = 1752.9032
$ws.Range("J132").Value = 61138.176
$ws.Range("K132").Value = 5258.7096
$ws.Range("L132").Value = 183414.528
$ws.Range("M132").Value = -2728.7096
$ws.Range("N132").Value = -188474.528

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 107.666664
$ws.Range("I23").Value = 34.4
$ws.Range("J23").Value = 144.3
$ws.Range("K23").Value = 103.2
$ws.Range("L23").Value = 432.9
$ws.Range("M23").Value = 131.8
$ws.Range("N23").Value = -902.9000000000001
$ws.Range("H132").Value = 1814.1818
$ws.Range("I132").Value = 1569.625
$ws.Range("J132").Value = 2466.3333
$ws.Range("K132").Value = 14126.625
$ws.Range("L132").Value = 22196.9997
$ws.Range("M132").Value = -11596.625
$ws.Range("N132").Value = -27256.9997
$ws.Range("H133").Value = 4205.8184
$ws.Range("I133").Value = 3090
$ws.Range("K133").Value = 9270
$ws.Range("M133").Value = -4210
$ws.Range("H140").Value = 2163.8643
$ws.Range("I140").Value = 1976.4706
$ws.Range("J140").Value = 2299.4255
$ws.Range("K140").Value = 5929.4118
$ws.Range("L140").Value = 6898.2765
$ws.Range("M140").Value = -749.4117999999999
$ws.Range("N140").Value = -17258.2765
$ws.Range("H141").Value = 7254.706
$ws.Range("I141").Value = 3652.5
$ws.Range("K141").Value = 10957.5
$ws.Range("M141").Value = -5777.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 7500
$ws.Range("I21").Value = 5000
$ws.Range("J21").Value = 10000
$ws.Range("K21").Value = 5000
$ws.Range("L21").Value = 10000
$ws.Range("M21").Value = -4827
$ws.Range("N21").Value = -10346
$ws.Range("H30").Value = 7500
$ws.Range("I30").Value = 5000
$ws.Range("J30").Value = 10000
$ws.Range("K30").Value = 5000
$ws.Range("L30").Value = 10000
$ws.Range("M30").Value = -4895
$ws.Range("N30").Value = -10210

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 23070.834
$ws.Range("J63").Value = 23070.834
$ws.Range("L63").Value = 23070.834
$ws.Range("N63").Value = -24568.834
$ws.Range("H66").Value = 23070.834
$ws.Range("J66").Value = 23070.834
$ws.Range("L66").Value = 69212.50199999999
$ws.Range("N66").Value = -76700.50199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 125688.125
$ws.Range("I100").Value = 84083.664
$ws.Range("J100").Value = 250501.5
$ws.Range("K100").Value = 168167.328
$ws.Range("L100").Value = 501003
$ws.Range("M100").Value = -167626.328
$ws.Range("N100").Value = -502085
$ws.Range("H113").Value = 1707
$ws.Range("I113").Value = 1264
$ws.Range("J113").Value = 1873.125
$ws.Range("K113").Value = 3792
$ws.Range("L113").Value = 5619.375
$ws.Range("M113").Value = -1622
$ws.Range("N113").Value = -9959.375
$ws.Range("H136").Value = 33813.85
$ws.Range("I136").Value = 20139.27
$ws.Range("J136").Value = 112822.555
$ws.Range("K136").Value = 60417.81
$ws.Range("L136").Value = 338467.665
$ws.Range("M136").Value = -57867.81
$ws.Range("N136").Value = -343567.665
